# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending order of "Periodo Mora" values for rows 16..26 (was descending)
$periods = @("1701","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# The "Valor Mora" amounts for row 16 and row 26 swap along with the period reorder
$ws.Range("F16").Value = 27578
$ws.Range("F26").Value = 29509
